$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 7: Inscritos (E7) 26 -> 27
$ws.Range("E7").Value = 27

# Row 8: Inscritos (E8) 39 -> 40, Pagos (F8) 13 -> 14, Inscricoes homologadas (H8) 13 -> 14
$ws.Range("E8").Value = 40
$ws.Range("F8").Value = 14
$ws.Range("H8").Value = 14

# Row 12: Inscritos (E12) 25 -> 26, Pagos (F12) 9 -> 10, Inscricoes homologadas (H12) 9 -> 10
$ws.Range("E12").Value = 26
$ws.Range("F12").Value = 10
$ws.Range("H12").Value = 10

# Row 14: Inscritos (E14) 35 -> 36
$ws.Range("E14").Value = 36

$wb.Save()
